# Update simulation result values on Sheet1
# - correction damping - costfunctioninputs: areaUsage - design param: moduleRowSpacing

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = [double]"30"
$ws.Range("B4").Value = [double]"18732.19023536615"
$ws.Range("B5").Value = [double]"1485.548159853576"
$ws.Range("B6").Value = [double]"0.1383952194513708"
$ws.Range("B7").Value = [double]"3.090678116214709e-18"
$ws.Range("B8").Value = [double]"20217.60000000103"
$ws.Range("B9").Value = [double]"5664.220747500913"
$ws.Range("B10").Value = [double]"292204.1456634886"
$ws.Range("B11").Value = [double]"1.388618040753951"
$ws.Range("B12").Value = [double]"1.210034551702396e+20"
$ws.Range("B13").Value = [double]"0.3499999999999936"
$ws.Range("B14").Value = [double]"0.9999068389546928"
$ws.Range("B15").Value = [double]"0.0734706194925014"
